$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix wording / typo in microphone signal note (row 8)
$ws.Range("C8").Value = 'Подключение сигнала микрофона к AGND при отсутствии штекера'

# Add "Конфликтует с" cross references for R1 / R7 pull-up rows (rows 9 and 10)
$ws.Range("D9").Value = 'R7'
$ws.Range("D10").Value = 'R1'

# Mark R7 pull-up to 3.3V (MCU) as no longer used (row 10)
$ws.Range("C10").Value = 'Pull-up кнопок к 3.3V (MCU) - не используется'

# Fix typo "Подключенbе" -> "Подключение" in battery "-" to GND note (row 23)
$ws.Range("C23").Value = 'Подключение "-" аккумулятора к GND в отсутстсвие VT4'

$wb.Save()
